# Update "想去人数" (want-to-go count) figures by +1 on both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1271
$ws1.Range("F4").Value = 2766
$ws1.Range("F5").Value = 247

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1271
$ws4.Range("F6").Value = 2766
$ws4.Range("F8").Value = 247
